# Refresh the "paises" COVID dashboard (sheet "Pais") with the 28-Sep-2020
# 16:37 data pull: updated case counts for a batch of countries, plus three
# countries (Birmania, Principado de Andorra, Santa Lucia) that overtook
# their neighbours in the ranking and so swap places with them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 16:37"

# --- Row-by-row updates: label reorders + refreshed case counts ---
# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 7322426
$ws.Range("C4").Value = 1083
$ws.Range("D4").Value = 4560771
$ws.Range("E4").Value = 2552182
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 209473

# Row 5: 'India' -> 'India'
$ws.Range("B5").Value = 6083888
$ws.Range("C5").Value = 10540
$ws.Range("D5").Value = 5020881
$ws.Range("E5").Value = 967381
$ws.Range("G5").Value = 52
$ws.Range("H5").Value = 95626

# Row 25: 'Alemania' -> 'Alemania'
$ws.Range("B25").Value = 287269
$ws.Range("C25").Value = 931
$ws.Range("E25").Value = 26934

# Row 52: 'Portugal' -> 'Portugal'
$ws.Range("B52").Value = 74029
$ws.Range("C52").Value = 425
$ws.Range("D52").Value = 47884
$ws.Range("E52").Value = 24188
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 1957

# Row 61: 'Suiza' -> 'Suiza'
$ws.Range("E61").Value = 7881
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 2065

# Row 75: 'Serbia' -> 'Serbia'
$ws.Range("B75").Value = 33414
$ws.Range("C75").Value = 30
$ws.Range("E75").Value = 1130
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 748

# Row 80: 'Bosnia y Herzegovina' -> 'Bosnia y Herzegovina'
$ws.Range("B80").Value = 27001
$ws.Range("C80").Value = 81
$ws.Range("D80").Value = 20005
$ws.Range("E80").Value = 6167
$ws.Range("G80").Value = 7
$ws.Range("H80").Value = 829

# Row 93: 'Noruega' -> 'Noruega'
$ws.Range("B93").Value = 13741
$ws.Range("C93").Value = 43
$ws.Range("E93").Value = 2277
$ws.Range("G93").Value = 4
$ws.Range("H93").Value = 274

# Row 96: 'Malasia' -> 'Birmania'
$ws.Range("A96").Value = "Birmania"
$ws.Range("B96").Value = 11631
$ws.Range("C96").Value = 897
$ws.Range("D96").Value = 3073
$ws.Range("E96").Value = 8302
$ws.Range("G96").Value = 30
$ws.Range("H96").Value = 256

# Row 97: 'Namibia' -> 'Malasia'
$ws.Range("A97").Value = "Malasia"
$ws.Range("B97").Value = 11034
$ws.Range("C97").Value = 115
$ws.Range("D97").Value = 9889
$ws.Range("E97").Value = 1011
$ws.Range("H97").Value = 134

# Row 98: 'Birmania' -> 'Namibia'
$ws.Range("A98").Value = "Namibia"
$ws.Range("B98").Value = 11033
$ws.Range("D98").Value = 8776
$ws.Range("E98").Value = 2137
$ws.Range("H98").Value = 120

# Row 133: 'Trinidad yTobago' -> 'Trinidad yTobago'
$ws.Range("B133").Value = 4382
$ws.Range("C133").Value = 20
$ws.Range("D133").Value = 2275
$ws.Range("E133").Value = 2035
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 72

# Row 156: 'Belice' -> 'Principado de Andorra'
$ws.Range("A156").Value = "Principado de Andorra"
$ws.Range("B156").Value = 1966
$ws.Range("C156").Value = 130
$ws.Range("D156").Value = 1265
$ws.Range("E156").Value = 648
$ws.Range("H156").Value = 53

# Row 157: 'Principado de Andorra' -> 'Belice'
$ws.Range("A157").Value = "Belice"
$ws.Range("B157").Value = 1854
$ws.Range("C157").Value = 29
$ws.Range("D157").Value = 1196
$ws.Range("E157").Value = 634
$ws.Range("H157").Value = 24

# Row 163: 'Lesoto' -> 'Lesoto'
$ws.Range("B163").Value = 1565
$ws.Range("C163").Value = 7
$ws.Range("D163").Value = 822
$ws.Range("E163").Value = 708

# Row 164: 'Liberia' -> 'Liberia'
$ws.Range("B164").Value = 1342
$ws.Range("C164").Value = 3
$ws.Range("E164").Value = 39

# Row 184: 'Curazao' -> 'Curazao'
$ws.Range("B184").Value = 366
$ws.Range("C184").Value = 6
$ws.Range("D184").Value = 148
$ws.Range("E184").Value = 217

# Row 207: 'Timor Oriental' -> 'Santa Lucia'
$ws.Range("A207").Value = "Santa Lucia"

# Row 208: 'Santa Lucia' -> 'Timor Oriental'
$ws.Range("A208").Value = "Timor Oriental"
